$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.154.45"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").Value = "2.424.38"
$ws.Range("E3").Value = "  +3.13%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'554.78"
$ws.Range("E5").Value = "  +2.49%  "
$ws.Range("D6").Value = "'138.75"
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("E9").Value = "  +3.27%  "
$ws.Range("E10").Value = "  +4.32%  "
$ws.Range("E11").Value = "  +0.57%  "
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").Value = "2.854.93"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "60.066.14"
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "2.421.79"
$ws.Range("E17").Value = "  +3.22%  "
$ws.Range("D18").Value = "'11.40"
$ws.Range("E18").Value = "  +6.40%  "
$ws.Range("E19").Value = "  +2.22%  "
$ws.Range("D20").Value = "'332.76"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'65.17"
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("D25").Value = "'8.60"
$ws.Range("E25").Value = "  +2.91%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  +7.22%  "
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "'6.30"
$ws.Range("E30").Value = "  +2.68%  "
$ws.Range("D31").Value = "'168.78"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("D33").Value = "'18.72"
$ws.Range("E33").Value = "  +1.80%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +5.43%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +11.81%  "
$ws.Range("D40").Value = "'321.53"
$ws.Range("E40").Value = "  +11.16%  "
$ws.Range("D41").Value = "'39.48"
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "'140.38"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "'0.0962"
$ws.Range("E44").Value = "  +1.24%  "
$ws.Range("D45").Value = "'0.0522"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("D46").Value = "'19.58"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'0.410"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "'17.80"
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("E51").Value = "  -0.19%  "
